$d = $word.ActiveDocument

# 1) Merge "So: ..../TB" + "-${soVB}" into a single run "So: ..../TB-${soVB}"
$d.Content.Find.Execute('Số: …../TB-${soVB}', $true, $false, $false, $false, $false, $true, 1, $false, 'Số: …../TB-${soVB}', 2)

# 2) Add the year "2018" to the date line
$d.Content.Find.Execute('., ngày      tháng      năm', $true, $false, $false, $false, $false, $true, 1, $false, '., ngày      tháng      năm 2018', 2)

# 3) Insert "dia chi : ${diaChi}," before "voi noi dung :"
$d.Content.Find.Execute('Ngày ${ngayTiepNhan}, ${coQuanTiepNhan} đã nhận được đơn tố cáo của ${nguoiToCao} với nội dung : ', $true, $false, $false, $false, $false, $true, 1, $false, 'Ngày ${ngayTiepNhan}, ${coQuanTiepNhan} đã nhận được đơn tố cáo của ${nguoiToCao}, địa chỉ : ${diaChi}, với nội dung : ', 2)

# 4) Replace the blank "Ly do: ....." placeholder with the template variable
$d.Content.Find.Execute('Lý do: …………………………………………………………………………………………..(2)', $true, $false, $false, $false, $false, $true, 1, $false, 'Lý do: ${lyDoDinhChi}.', 2)

# 5) Remove the trailing footnote (2) paragraphs that explained the now-removed placeholder
$delStart = $d.Content
$delStart.Find.Execute('(2) Căn cứ pháp lý đã áp dụng để không thụ lý giải quyết tố cáo (điều, khoản, tên văn bản và', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$delRange = $d.Range($delStart.Start, $d.Content.End)
$delRange.Delete()

# 6) Nudge the small underline drawing below "THONG BAO" back in line with the new text metrics
$shp = $d.Shapes(1)
$shp.Width = 64.65
$shp.Height = 0.65
